$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "nityaranjn55623@gmail.com"
$ws.Range("C3").Value = "abh09082444@gmail.com"
